# This script reproduces the numeric updates to the per-job "Profits" sheets
# (columns H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
# K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ) coming from a
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

# --- Sheet "ALC" (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(74, 8).Value = 10203.692  # H74
$ws.Cells.Item(74, 10).Value = 19166.666  # J74
$ws.Cells.Item(74, 12).Value = 19166.666  # L74
$ws.Cells.Item(74, 14).Value = -21038.666  # N74
$ws.Cells.Item(77, 8).Value = 10203.692  # H77
$ws.Cells.Item(77, 10).Value = 19166.666  # J77
$ws.Cells.Item(77, 12).Value = 95833.33  # L77
$ws.Cells.Item(77, 14).Value = -105193.33  # N77
$ws.Cells.Item(100, 8).Value = 7365.08  # H100
$ws.Cells.Item(100, 9).Value = 2383  # I100
$ws.Cells.Item(100, 10).Value = 13705.909  # J100
$ws.Cells.Item(100, 11).Value = 2383  # K100
$ws.Cells.Item(100, 12).Value = 13705.909  # L100
$ws.Cells.Item(100, 13).Value = -1842  # M100
$ws.Cells.Item(100, 14).Value = -14787.909  # N100
$ws.Cells.Item(127, 8).Value = 29520.25  # H127
$ws.Cells.Item(127, 9).Value = 37693.668  # I127
$ws.Cells.Item(127, 11).Value = 113081.004  # K127
$ws.Cells.Item(127, 13).Value = -108121.004  # M127
$ws.Cells.Item(131, 8).Value = 4805.467  # H131
$ws.Cells.Item(131, 9).Value = 2281.6667  # I131
$ws.Cells.Item(131, 11).Value = 6845.000100000001  # K131
$ws.Cells.Item(131, 13).Value = -1805.000100000001  # M131
$ws.Cells.Item(132, 8).Value = 2353.1538  # H132
$ws.Cells.Item(132, 9).Value = 1325.5667  # I132
$ws.Cells.Item(132, 11).Value = 3976.7001  # K132
$ws.Cells.Item(132, 13).Value = -1446.7001  # M132
$ws.Cells.Item(135, 8).Value = 976.6875  # H135
$ws.Cells.Item(135, 9).Value = 624.46344  # I135
$ws.Cells.Item(135, 10).Value = 3039.7144  # J135
$ws.Cells.Item(135, 11).Value = 5620.170959999999  # K135
$ws.Cells.Item(135, 12).Value = 27357.4296  # L135
$ws.Cells.Item(135, 13).Value = -3085.170959999999  # M135
$ws.Cells.Item(135, 14).Value = -32427.4296  # N135
$ws.Cells.Item(138, 8).Value = 3778.922  # H138
$ws.Cells.Item(138, 10).Value = 5646.5415  # J138
$ws.Cells.Item(138, 12).Value = 16939.6245  # L138
$ws.Cells.Item(138, 14).Value = -27219.6245  # N138
$ws.Cells.Item(141, 8).Value = 1610.9524  # H141
$ws.Cells.Item(141, 9).Value = 1610.9524  # I141
$ws.Cells.Item(141, 11).Value = 4832.857199999999  # K141
$ws.Cells.Item(141, 13).Value = 347.1428000000005  # M141

# --- Sheet "ARM" (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(17, 8).Value = 10000  # H17
$ws.Cells.Item(17, 10).Value = 10000  # J17
$ws.Cells.Item(17, 12).Value = 10000  # L17
$ws.Cells.Item(17, 14).Value = -10346  # N17
$ws.Cells.Item(45, 8).Value = 8087.909  # H45
$ws.Cells.Item(45, 9).Value = 5995.2856  # I45
$ws.Cells.Item(45, 10).Value = 11750  # J45
$ws.Cells.Item(45, 11).Value = 5995.2856  # K45
$ws.Cells.Item(45, 12).Value = 11750  # L45
$ws.Cells.Item(45, 13).Value = -5618.2856  # M45
$ws.Cells.Item(45, 14).Value = -12504  # N45
$ws.Cells.Item(132, 8).Value = 2584.3374  # H132
$ws.Cells.Item(132, 9).Value = 863.92645  # I132
$ws.Cells.Item(132, 11).Value = 2591.77935  # K132
$ws.Cells.Item(132, 13).Value = -61.77935000000025  # M132

# --- Sheet "BSM" (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(35, 8).Value = 102583.164  # H35
$ws.Cells.Item(35, 10).Value = 102583.164  # J35
$ws.Cells.Item(35, 12).Value = 102583.164  # L35
$ws.Cells.Item(35, 14).Value = -103203.164  # N35
$ws.Cells.Item(86, 8).Value = 896299.75  # H86
$ws.Cells.Item(86, 9).Value = 2126646.5  # I86
$ws.Cells.Item(86, 10).Value = 1502.1818  # J86
$ws.Cells.Item(86, 11).Value = 2126646.5  # K86
$ws.Cells.Item(86, 12).Value = 1502.1818  # L86
$ws.Cells.Item(86, 13).Value = -2125523.5  # M86
$ws.Cells.Item(86, 14).Value = -3748.1818  # N86
$ws.Cells.Item(89, 8).Value = 896299.75  # H89
$ws.Cells.Item(89, 9).Value = 2126646.5  # I89
$ws.Cells.Item(89, 10).Value = 1502.1818  # J89
$ws.Cells.Item(89, 11).Value = 10633232.5  # K89
$ws.Cells.Item(89, 12).Value = 7510.909000000001  # L89
$ws.Cells.Item(89, 13).Value = -10627616.5  # M89
$ws.Cells.Item(89, 14).Value = -18742.909  # N89
$ws.Cells.Item(134, 8).Value = 26787.697  # H134
$ws.Cells.Item(134, 9).Value = 2263.8064  # I134
$ws.Cells.Item(134, 10).Value = 90141.086  # J134
$ws.Cells.Item(134, 11).Value = 6791.4192  # K134
$ws.Cells.Item(134, 12).Value = 270423.258  # L134
$ws.Cells.Item(134, 13).Value = -4256.4192  # M134
$ws.Cells.Item(134, 14).Value = -275493.258  # N134

# --- Sheet "CRP" (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 676  # H22
$ws.Cells.Item(22, 9).Value = 650  # I22
$ws.Cells.Item(22, 10).Value = 702  # J22
$ws.Cells.Item(22, 11).Value = 650  # K22
$ws.Cells.Item(22, 12).Value = 702  # L22
$ws.Cells.Item(22, 13).Value = -300  # M22
$ws.Cells.Item(22, 14).Value = -1402  # N22
$ws.Cells.Item(31, 8).Value = 669992.9399999999  # H31
$ws.Cells.Item(31, 9).Value = 835241.2  # I31
$ws.Cells.Item(31, 10).Value = 9000  # J31
$ws.Cells.Item(31, 11).Value = 835241.2  # K31
$ws.Cells.Item(31, 12).Value = 9000  # L31
$ws.Cells.Item(31, 13).Value = -834946.2  # M31
$ws.Cells.Item(31, 14).Value = -9590  # N31
$ws.Cells.Item(34, 8).Value = 669992.9399999999  # H34
$ws.Cells.Item(34, 9).Value = 835241.2  # I34
$ws.Cells.Item(34, 10).Value = 9000  # J34
$ws.Cells.Item(34, 11).Value = 835241.2  # K34
$ws.Cells.Item(34, 12).Value = 9000  # L34
$ws.Cells.Item(34, 13).Value = -835039.2  # M34
$ws.Cells.Item(34, 14).Value = -9404  # N34
$ws.Cells.Item(99, 8).Value = 6771.643  # H99
$ws.Cells.Item(99, 9).Value = 7022.2856  # I99
$ws.Cells.Item(99, 11).Value = 7022.2856  # K99
$ws.Cells.Item(99, 13).Value = -5524.2856  # M99
$ws.Cells.Item(126, 8).Value = 6771.643  # H126
$ws.Cells.Item(126, 9).Value = 7022.2856  # I126
$ws.Cells.Item(126, 11).Value = 21066.8568  # K126
$ws.Cells.Item(126, 13).Value = -18596.8568  # M126
$ws.Cells.Item(132, 8).Value = 3115.4443  # H132
$ws.Cells.Item(132, 9).Value = 2103.3171  # I132
$ws.Cells.Item(132, 11).Value = 6309.951300000001  # K132
$ws.Cells.Item(132, 13).Value = -3779.951300000001  # M132

# --- Sheet "CUL" (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(51, 8).Value = 855  # H51
$ws.Cells.Item(51, 9).Value = 838.3333  # I51
$ws.Cells.Item(51, 11).Value = 2514.9999  # K51
$ws.Cells.Item(51, 13).Value = -2054.9999  # M51
$ws.Cells.Item(131, 8).Value = 4534.1577  # H131
$ws.Cells.Item(131, 10).Value = 6184.615  # J131
$ws.Cells.Item(131, 12).Value = 18553.845  # L131
$ws.Cells.Item(131, 14).Value = -28633.845  # N131
$ws.Cells.Item(132, 8).Value = 3128.1428  # H132
$ws.Cells.Item(132, 10).Value = 5783.1665  # J132
$ws.Cells.Item(132, 12).Value = 52048.4985  # L132
$ws.Cells.Item(132, 14).Value = -57108.4985  # N132

# --- Sheet "GSM" (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(21, 8).Value = 700000000  # H21
$ws.Cells.Item(21, 9).Value = 0  # I21
$ws.Cells.Item(21, 10).Value = 700000000  # J21
$ws.Cells.Item(21, 11).Value = 0  # K21
$ws.Cells.Item(21, 12).Value = 700000000  # L21
$ws.Cells.Item(21, 13).ClearContents()  # M21 removed
$ws.Cells.Item(21, 14).Value = -700000346  # N21
$ws.Cells.Item(30, 8).Value = 700000000  # H30
$ws.Cells.Item(30, 9).Value = 0  # I30
$ws.Cells.Item(30, 10).Value = 700000000  # J30
$ws.Cells.Item(30, 11).Value = 0  # K30
$ws.Cells.Item(30, 12).Value = 700000000  # L30
$ws.Cells.Item(30, 13).ClearContents()  # M30 removed
$ws.Cells.Item(30, 14).Value = -700000210  # N30
$ws.Cells.Item(122, 8).Value = 3609.2307  # H122
$ws.Cells.Item(122, 9).Value = 3525.4443  # I122
$ws.Cells.Item(122, 10).Value = 3797.75  # J122
$ws.Cells.Item(122, 11).Value = 10576.3329  # K122
$ws.Cells.Item(122, 12).Value = 11393.25  # L122
$ws.Cells.Item(122, 13).Value = -8126.332900000001  # M122
$ws.Cells.Item(122, 14).Value = -16293.25  # N122

# --- Sheet "LTW" (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 483210.06  # H7
$ws.Cells.Item(7, 9).Value = 7409.0835  # I7
$ws.Cells.Item(7, 11).Value = 7409.0835  # K7
$ws.Cells.Item(7, 13).Value = -7297.0835  # M7
$ws.Cells.Item(40, 8).Value = 4734.467  # H40
$ws.Cells.Item(40, 9).Value = 2889.625  # I40
$ws.Cells.Item(40, 10).Value = 6842.857  # J40
$ws.Cells.Item(40, 11).Value = 2889.625  # K40
$ws.Cells.Item(40, 12).Value = 6842.857  # L40
$ws.Cells.Item(40, 13).Value = -2753.625  # M40
$ws.Cells.Item(40, 14).Value = -7114.857  # N40
$ws.Cells.Item(46, 8).Value = 4781.4736  # H46
$ws.Cells.Item(46, 9).Value = 4259  # I46
$ws.Cells.Item(46, 10).Value = 5499.875  # J46
$ws.Cells.Item(46, 11).Value = 4259  # K46
$ws.Cells.Item(46, 12).Value = 5499.875  # L46
$ws.Cells.Item(46, 13).Value = -4071  # M46
$ws.Cells.Item(46, 14).Value = -5875.875  # N46
$ws.Cells.Item(55, 8).Value = 976  # H55
$ws.Cells.Item(55, 9).Value = 585.75  # I55
$ws.Cells.Item(55, 10).Value = 1669.7778  # J55
$ws.Cells.Item(55, 11).Value = 585.75  # K55
$ws.Cells.Item(55, 12).Value = 1669.7778  # L55
$ws.Cells.Item(55, 13).Value = -412.75  # M55
$ws.Cells.Item(55, 14).Value = -2015.7778  # N55
$ws.Cells.Item(93, 8).Value = 2163.1667  # H93
$ws.Cells.Item(93, 9).Value = 2105.8  # I93
$ws.Cells.Item(93, 11).Value = 2105.8  # K93
$ws.Cells.Item(93, 13).Value = -857.8000000000002  # M93
$ws.Cells.Item(126, 8).Value = 483210.06  # H126
$ws.Cells.Item(126, 9).Value = 7409.0835  # I126
$ws.Cells.Item(126, 11).Value = 22227.2505  # K126
$ws.Cells.Item(126, 13).Value = -19757.2505  # M126
$ws.Cells.Item(136, 8).Value = 6405.4443  # H136
$ws.Cells.Item(136, 9).Value = 3386  # I136
$ws.Cells.Item(136, 10).Value = 7268.143  # J136
$ws.Cells.Item(136, 11).Value = 10158  # K136
$ws.Cells.Item(136, 12).Value = 21804.429  # L136
$ws.Cells.Item(136, 13).Value = -7608  # M136
$ws.Cells.Item(136, 14).Value = -26904.429  # N136

# --- Sheet "WVR" (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(64, 8).Value = 78019.2  # H64
$ws.Cells.Item(67, 8).Value = 78019.2  # H67
$ws.Cells.Item(126, 8).Value = 3600.353  # H126
$ws.Cells.Item(126, 9).Value = 2070.7  # I126
$ws.Cells.Item(126, 11).Value = 6212.099999999999  # K126
$ws.Cells.Item(126, 13).Value = -3742.099999999999  # M126
$ws.Cells.Item(132, 8).Value = 27900.846  # H132
$ws.Cells.Item(132, 9).Value = 1207.7667  # I132
$ws.Cells.Item(132, 11).Value = 3623.300099999999  # K132
$ws.Cells.Item(132, 13).Value = -1093.300099999999  # M132
